$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.751841
$ws.Range("N2").Value = 3.503682
$ws.Range("O2").Value = 0.213779419532536
$ws.Range("P2").Value = 0.1765245924501479
$ws.Range("Q2").Value = 0.266568885765
$ws.Range("R2").Value = 1.59941331459
$ws.Range("S2").Value = 0.213779419532536
$ws.Range("T2").Value = 0.1765245924501479

# Row 3
$ws.Range("O3").Value = 0.2839067349822242
$ws.Range("P3").Value = 0.3516464830630844
$ws.Range("S3").Value = 0.2839067349822242
$ws.Range("T3").Value = 0.3516464830630844

# Row 4
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.188968
$ws.Range("N4").Value = 0.5669040000000001
$ws.Range("O4").Value = 0.02306000907058589
$ws.Range("P4").Value = 0.02856209483576382
$ws.Range("Q4").Value = 0.02875431572
$ws.Range("R4").Value = 0.25878884148
$ws.Range("S4").Value = 0.02306000907058589
$ws.Range("T4").Value = 0.02856209483576382

# Row 5
$ws.Range("M5").Value = 2.9838935
$ws.Range("N5").Value = 5.967787
$ws.Range("O5").Value = 0.3641283771626004
$ws.Range("P5").Value = 0.3006725975714379
$ws.Range("Q5").Value = 0.4540441544275
$ws.Range("R5").Value = 2.724264926565
$ws.Range("S5").Value = 0.3641283771626004
$ws.Range("T5").Value = 0.3006725975714379

# Row 6
$ws.Range("M6").Value = 0.6503553333333333
$ws.Range("N6").Value = 1.951066
$ws.Range("O6").Value = 0.07936370118628855
$ws.Range("P6").Value = 0.0982997687841934
$ws.Range("Q6").Value = 0.09896131929666666
$ws.Range("R6").Value = 0.89065187367
$ws.Range("S6").Value = 0.07936370118628855
$ws.Range("T6").Value = 0.0982997687841934

# Row 7
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.293054
$ws.Range("N7").Value = 0.879162
$ws.Range("O7").Value = 0.03576175806576498
$ws.Range("P7").Value = 0.04429446329537239
$ws.Range("Q7").Value = 0.04459256190999999
$ws.Range("R7").Value = 0.40133305719
$ws.Range("S7").Value = 0.03576175806576498
$ws.Range("T7").Value = 0.04429446329537239
